# Daily attendance processing - 2026-01-24 21:01:33
# Swap the order of names in the "Recorded By" column (G) from
# "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# for all rows where that exact value is present.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$targetRows = @(2,3,4,5,6,7,16,17,22,23,37,38,43,44,58,59,64,65,79,80,85,86,87,88,89,90,99,100,105,106,107,108,109,110,119,120,125,126,127,128,129,130,139,140,145,146,147,148,149,150,159,160,165,166,167,168,169,170,179,180,185,186,200,201,206,207,221,222,227,228,242,243)

foreach ($row in $targetRows) {
    $cell = $ws.Range("G$row")
    if ($cell.Value2 -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
